$wb = $excel.ActiveWorkbook

# --- Sheet 1: By_Odds_Bin ---
$ws1 = $wb.Worksheets.Item("By_Odds_Bin")

# Row 2: (0, 5]
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 0
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 0
$ws1.Range("F2").Value = ""

# Row 3: (5, 10]
$ws1.Range("B3").Value = 32
$ws1.Range("C3").Value = 17
$ws1.Range("D3").Value = 43
$ws1.Range("E3").Value = -26
$ws1.Range("F3").Value = 18.8

# Row 4: (10, 15]
$ws1.Range("B4").Value = 8
$ws1.Range("C4").Value = 16
$ws1.Range("D4").Value = 22
$ws1.Range("E4").Value = -6
$ws1.Range("F4").Value = 25

# --- Sheet 2: By_Field_Size ---
$ws2 = $wb.Worksheets.Item("By_Field_Size")

# Row 2: 1-4
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Value = ""

# Row 3: 5
$ws2.Range("B3").Value = 0
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = ""

# Row 4: 6
$ws2.Range("B4").Value = 0
$ws2.Range("C4").Value = 0
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = ""

# Row 5: 7
$ws2.Range("B5").Value = 40
$ws2.Range("C5").Value = 33
$ws2.Range("D5").Value = 65
$ws2.Range("E5").Value = -32
$ws2.Range("F5").Value = 20

# --- Sheet 3: By_Track ---
$ws3 = $wb.Worksheets.Item("By_Track")

# Row 2: NEWMARKET -> NEWBURY
$ws3.Range("A2").Value = "NEWBURY"
$ws3.Range("B2").Value = 40
$ws3.Range("C2").Value = 33
$ws3.Range("D2").Value = 65
$ws3.Range("E2").Value = -32
$ws3.Range("F2").Value = 20
